$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting (avoid Excel auto-parsing numeric-looking strings
# in the Price/Volume columns as numbers) by forcing Text number format on
# every D/E cell we touch, then assigning the literal string value.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue "D2" "308.78"
Set-TextValue "E2" "0.29%"
Set-TextValue "D3" "41.07"
Set-TextValue "E3" "-1.52%"
Set-TextValue "E4" "2.05%"
Set-TextValue "E6" "-0.34%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D7" "0.9184"
Set-TextValue "E7" "1.53%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D8" "2.444"
Set-TextValue "E8" "-2.11%"
Set-TextValue "D9" "0.1238"
Set-TextValue "E9" "13.14%"
Set-TextValue "D10" "0.1839"
Set-TextValue "E10" "4.24%"
Set-TextValue "D11" "0.09176"
Set-TextValue "E11" "-0.42%"
Set-TextValue "D12" "0.04290"
Set-TextValue "E12" "0.42%"
Set-TextValue "D13" "0.1052"
Set-TextValue "E13" "0.05%"
Set-TextValue "D14" "0.001262"
Set-TextValue "E14" "0.21%"
Set-TextValue "D15" "0.005817"
Set-TextValue "E15" "-0.08%"
Set-TextValue "D17" "3.352"
Set-TextValue "E17" "-0.28%"
Set-TextValue "D18" "4.320"
Set-TextValue "E18" "1.65%"
Set-TextValue "E19" "1.22%"
Set-TextValue "D20" "7.216"
Set-TextValue "E20" "10.09%"
Set-TextValue "E21" "1.50%"
Set-TextValue "D22" "0.2893"
Set-TextValue "E22" "7.86%"
Set-TextValue "D23" "0.04078"
Set-TextValue "E23" "-2.18%"
Set-TextValue "E24" "3.64%"
Set-TextValue "D25" "0.004140"
Set-TextValue "E25" "1.05%"
Set-TextValue "E26" "-2.15%"
Set-TextValue "D38" "0.02449"
Set-TextValue "E38" "1.20%"
Set-TextValue "D39" "0.05298"
Set-TextValue "E39" "1.94%"
Set-TextValue "E40" "1.00%"
Set-TextValue "E41" "0.97%"
Set-TextValue "D42" "0.006822"
Set-TextValue "E42" "-1.90%"
Set-TextValue "E43" "-0.34%"
Set-TextValue "D44" "0.007789"
Set-TextValue "E44" "-3.46%"
Set-TextValue "E45" "0.52%"
Set-TextValue "E46" "-1.01%"
Set-TextValue "E47" "0.16%"
Set-TextValue "D48" "0.1699"
Set-TextValue "E48" "1,757.13%"
Set-TextValue "E49" "-2.43%"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.16%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.16%"
